# Update the "totaalstand" (standings) table on Sheet1 with the latest
# tournament results: scores/180s/100+ finishes recalculated, several
# players re-ranked/re-ordered, and two new players (Gijs Tromp, Nick
# Fitzpatrick) added at the bottom (rows 19-20), extending the used range
# from A1:F18 to A1:F20.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 'Nathan May'
$ws.Range("C2").Value = 33
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 38
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 'Yannick den Daggelder'
$ws.Range("C3").Value = 33
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 34
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 'Louis Tweddle'
$ws.Range("C4").Value = 17
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 18
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 'Alessandro Delia'
$ws.Range("C5").Value = 17
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 17
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 'Remco Riem'
$ws.Range("C6").Value = 16
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 17
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 'Rocky Van Den Eeckhoudt'
$ws.Range("C7").Value = 12
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 14
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 'Milan Schoenmakers'
$ws.Range("C8").Value = 11
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 12
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 'Afendi Kelana'
$ws.Range("C9").Value = 8
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 11
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 'Diego Meerveld'
$ws.Range("C10").Value = 9
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 9
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 'Francesco Petruccelli'
$ws.Range("C11").Value = 8
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 9
$ws.Range("A12").Value = 9
$ws.Range("B12").Value = 'martin Berry'
$ws.Range("C12").Value = 8
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 9
$ws.Range("A13").Value = 9
$ws.Range("B13").Value = 'Nick Fitzpatrick'
$ws.Range("C13").Value = 8
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 9
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 'Dartin Dan'
$ws.Range("C14").Value = 6
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 7
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 'Magnus Gladh'
$ws.Range("C15").Value = 4
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 4
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 'Max Walter'
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 3
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = 'Robin Willis'
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 3
$ws.Range("A18").Value = 15
$ws.Range("B18").Value = 'Tristan Snoep'
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 3
$ws.Range("A19").Value = 15
$ws.Range("B19").Value = 'Gijs Tromp'
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 3
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = 'Lee Dolan'
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 1
